$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 in package order, 1st worksheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6771
$ws1.Range("F3").Value = 401
$ws1.Range("F4").Value = 111
$ws1.Range("F5").Value = 158
$ws1.Range("F8").Value = 592

# Sheet "全部类型" (sheet4 in package order, 4th worksheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6772
$ws4.Range("G2").Value = 62
$ws4.Range("F3").Value = 401
$ws4.Range("F5").Value = 111
$ws4.Range("F6").Value = 158
$ws4.Range("F10").Value = 592
